$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 235, shifting existing rows 235-292 down to 236-293
$ws.Rows(235).Insert()

# Populate the newly inserted row 235 with the new data record
$ws.Cells.Item(235, 1).Value = 7
$ws.Cells.Item(235, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(235, 3).Value = "Ñuble"
$ws.Cells.Item(235, 4).Value = (Get-Date -Year 2023 -Month 3 -Day 20 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(235, 5).Value = 16
$ws.Cells.Item(235, 6).Value = 100112032
$ws.Cells.Item(235, 7).Value = "Zapallo italiano"
$ws.Cells.Item(235, 8).Value = "Sin especificar"
$ws.Cells.Item(235, 9).Value = "Primera"
$ws.Cells.Item(235, 10).Value = 60
$ws.Cells.Item(235, 11).Value = 6000
$ws.Cells.Item(235, 12).Value = 6000
$ws.Cells.Item(235, 13).Value = 6000
$ws.Cells.Item(235, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(235, 15).Value = "Región del Maule"
$ws.Cells.Item(235, 16).Value = 120
$ws.Cells.Item(235, 17).Value = 50
$ws.Cells.Item(235, 18).Value = "Hortaliza"
